# "modify show batch img" - add an extra batch-image score column on the
# results sheet (Trang_tinh2 / "Sheet2"): insert a new column before F,
# fill in its first value, and move the selection to the newly inserted
# area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trang_tính2")

# Insert a new column before column F - this shifts the old F:I columns
# (and their column-width definitions) one slot to the right (G:J), and
# carries formatting from the neighboring column E down each row
# (e.g. row 8's styled-but-empty cell).
$ws.Columns("F").Insert()

# New data point for row 6 that now lives in the freshly inserted column F.
$ws.Range("F6").Value = 0.8163

# Leave the selection on the newly added cell.
$ws.Range("F7").Select() | Out-Null
